$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2.3
$ws.Range("L2").Value = 3.1
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.36
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.67
$ws.Range("T2").Value = 1.22
$ws.Range("Y2").Value = 8.5
$ws.Range("AN2").Value = 21
